$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title: the heading was previously split across two runs
#    ("...Analysis" / " using HR Data"). Re-typing the full phrase over
#    both runs collapses them into a single run with identical text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Project Summary – Employee Attrition Analysis using HR Data",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Project Summary – Employee Attrition Analysis using HR Data", 2)

# ---------------------------------------------------------------------
# 2. "Data Cleaning:" bullet: merge the three runs into one.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Data Cleaning: Dropped irrelevant columns (like EmployeeID).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data Cleaning: Dropped irrelevant columns (like EmployeeID).", 2)

# ---------------------------------------------------------------------
# 3. "Extracted insights..." bullet: merge the three runs into one.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Extracted insights to link key factors with attrition behaviour.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Extracted insights to link key factors with attrition behaviour.", 2)

# ---------------------------------------------------------------------
# 4. Conclusion bullet 1: insert " made" after "The Decision Tree model"
#    as its own run, splitting what was one run into three.
# ---------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute("The Decision Tree model")
$afterModel = $rng4.Duplicate
$afterModel.Collapse(0)
$afterModel.InsertAfter(" made")

# Force the newly-typed " made" text to live in its own run element
# (same text format as its neighbours, just a distinct <w:r>) by
# toggling a character property on/off across exactly that span.
$rngMade = $d.Content
$rngMade.Find.Execute(" made")
$madeOnly = $rngMade.Duplicate
$madeOnly.Font.Bold = $true
$madeOnly.Font.Bold = $false

# ---------------------------------------------------------------------
# 5. Conclusion bullet 3: extend "risk." with an explanatory
#    parenthetical, ending up as three runs: the lead-in text, the new
#    "risk (...)" phrase, and the trailing period.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Zone A and female employees show slightly higher risk.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zone A and female employees show slightly higher risk (i.e.; Employees who have a lower work life balance score and a lower rating and engagement score, they tend to leave on their own).", 2)

$rng5 = $d.Content
$rng5.Find.Execute("risk (i.e.; Employees who have a lower work life balance score and a lower rating and engagement score, they tend to leave on their own)")
$riskOnly = $rng5.Duplicate
$riskOnly.Font.Bold = $true
$riskOnly.Font.Bold = $false
